$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 49.6
$ws.Range("I8").Value = 32.88889
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 98.66667000000001
$ws.Range("L8").Value = 600
$ws.Range("M8").Value = 40.33332999999999
$ws.Range("N8").Value = -878

$ws.Range("H42").Value = 52728.316
$ws.Range("J42").Value = 100118
$ws.Range("L42").Value = 300354
$ws.Range("N42").Value = -300814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 200
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H5").Value = 276
$ws.Range("I5").Value = 276
$ws.Range("K5").Value = 276
$ws.Range("M5").Value = -164

$ws.Range("H12").Value = 203
$ws.Range("I12").Value = 203
$ws.Range("K12").Value = 203
$ws.Range("M12").Value = -30

$ws.Range("H17").Value = 70009
$ws.Range("J17").Value = 70009
$ws.Range("L17").Value = 70009
$ws.Range("N17").Value = -70355

$ws.Range("H32").Value = 1695149.4
$ws.Range("I32").Value = 2055193.6
$ws.Range("J32").Value = 5711.077
$ws.Range("K32").Value = 2055193.6
$ws.Range("L32").Value = 5711.077
$ws.Range("M32").Value = -2054906.6
$ws.Range("N32").Value = -6285.077

$ws.Range("H132").Value = 15815454
$ws.Range("I132").Value = 11288423
$ws.Range("J132").Value = 36510452
$ws.Range("K132").Value = 33865269
$ws.Range("L132").Value = 109531356
$ws.Range("M132").Value = -33862739
$ws.Range("N132").Value = -109536416

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 276
$ws.Range("I4").Value = 276
$ws.Range("K4").Value = 276
$ws.Range("M4").Value = -161

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H134").Value = 176473260
$ws.Range("I134").Value = 250004000
$ws.Range("J134").Value = 29411764
$ws.Range("K134").Value = 750012000
$ws.Range("L134").Value = 88235292
$ws.Range("M134").Value = -750009465
$ws.Range("N134").Value = -88240362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 53.81818
$ws.Range("I7").Value = 42
$ws.Range("J7").Value = 63.666668
$ws.Range("K7").Value = 42
$ws.Range("L7").Value = 63.666668
$ws.Range("M7").Value = 71
$ws.Range("N7").Value = -289.666668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 292
$ws.Range("I2").Value = 52.125
$ws.Range("K2").Value = 312.75
$ws.Range("M2").Value = -199.75

$ws.Range("H4").Value = 5375.75
$ws.Range("I4").Value = 750.5
$ws.Range("J4").Value = 10001
$ws.Range("K4").Value = 2251.5
$ws.Range("L4").Value = 30003
$ws.Range("M4").Value = -2139.5
$ws.Range("N4").Value = -30227

$ws.Range("H5").Value = 1751631.4
$ws.Range("I5").Value = 1479710.4
$ws.Range("K5").Value = 4439131.199999999
$ws.Range("M5").Value = -4439019.199999999

$ws.Range("H9").Value = 4428.5713
$ws.Range("I9").Value = 1500
$ws.Range("J9").Value = 4916.6665
$ws.Range("K9").Value = 4500
$ws.Range("L9").Value = 14749.9995
$ws.Range("M9").Value = -4276
$ws.Range("N9").Value = -15197.9995

$ws.Range("H10").Value = 494.85715
$ws.Range("I10").Value = 200
$ws.Range("K10").Value = 600
$ws.Range("M10").Value = -461

$ws.Range("H15").Value = 230
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 230
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 690
$ws.Range("N15").Value = -970
$ws.Range("M15").ClearContents()

$ws.Range("H16").Value = 900
$ws.Range("I16").Value = 900
$ws.Range("K16").Value = 2700
$ws.Range("M16").Value = -2527

$ws.Range("H19").Value = 35714684
$ws.Range("I19").Value = 35714684
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 107144052
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -107143878
$ws.Range("N19").ClearContents()

$ws.Range("H20").Value = 800
$ws.Range("I20").Value = 800
$ws.Range("K20").Value = 2400
$ws.Range("M20").Value = -2173

$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 3000
$ws.Range("N21").Value = -3346
$ws.Range("M21").ClearContents()

$ws.Range("H22").Value = 1157.8948
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1157.8948
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3473.6844
$ws.Range("N22").Value = -3811.6844
$ws.Range("M22").ClearContents()

$ws.Range("H26").Value = 156.5625
$ws.Range("I26").Value = 103.21429
$ws.Range("J26").Value = 530
$ws.Range("K26").Value = 309.64287
$ws.Range("L26").Value = 1590
$ws.Range("M26").Value = -21.64287000000002
$ws.Range("N26").Value = -2166

$ws.Range("H27").Value = 1157.8948
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1157.8948
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3473.6844
$ws.Range("N27").Value = -3677.6844
$ws.Range("M27").ClearContents()

$ws.Range("H32").Value = 1364.6666
$ws.Range("I32").Value = 1561.6
$ws.Range("K32").Value = 4684.799999999999
$ws.Range("M32").Value = -4401.799999999999

$ws.Range("H33").Value = 325
$ws.Range("I33").Value = 150
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 900
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = -617
$ws.Range("N33").Value = -3566

$ws.Range("H34").Value = 1930.5294
$ws.Range("I34").Value = 819.8333
$ws.Range("J34").Value = 2536.3635
$ws.Range("K34").Value = 2459.4999
$ws.Range("L34").Value = 7609.0905
$ws.Range("M34").Value = -2375.4999
$ws.Range("N34").Value = -7777.0905

$ws.Range("H38").Value = 20053.6
$ws.Range("I38").Value = 25059.5
$ws.Range("K38").Value = 75178.5
$ws.Range("M38").Value = -74831.5

$ws.Range("H39").Value = 2541.6667
$ws.Range("J39").Value = 2960
$ws.Range("L39").Value = 8880
$ws.Range("N39").Value = -9468

$ws.Range("H40").Value = 4778.5713
$ws.Range("I40").Value = 10510
$ws.Range("K40").Value = 42040
$ws.Range("M40").Value = -41971

$ws.Range("H44").Value = 1027.5
$ws.Range("I44").Value = 1000
$ws.Range("J44").Value = 1039.2858
$ws.Range("K44").Value = 3000
$ws.Range("L44").Value = 3117.8574
$ws.Range("M44").Value = -2602
$ws.Range("N44").Value = -3913.8574

$ws.Range("H46").Value = 1761.6296
$ws.Range("I46").Value = 1313
$ws.Range("J46").Value = 1863.591
$ws.Range("K46").Value = 3939
$ws.Range("L46").Value = 5590.772999999999
$ws.Range("M46").Value = -3848
$ws.Range("N46").Value = -5772.772999999999

$ws.Range("H50").Value = 225
$ws.Range("I50").Value = 118.75
$ws.Range("J50").Value = 650
$ws.Range("K50").Value = 356.25
$ws.Range("L50").Value = 1950
$ws.Range("M50").Value = 124.75
$ws.Range("N50").Value = -2912

$ws.Range("H51").Value = 2218.6365
$ws.Range("J51").Value = 2218.6365
$ws.Range("L51").Value = 6655.9095
$ws.Range("N51").Value = -7575.9095

$ws.Range("H53").Value = 225
$ws.Range("I53").Value = 118.75
$ws.Range("J53").Value = 650
$ws.Range("K53").Value = 356.25
$ws.Range("L53").Value = 1950
$ws.Range("M53").Value = 124.75
$ws.Range("N53").Value = -2912

$ws.Range("H54").Value = 5000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 5000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16118
$ws.Range("M54").ClearContents()

$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H58").Value = 2857.1428
$ws.Range("J58").Value = 3166.6667
$ws.Range("L58").Value = 9500.000100000001
$ws.Range("N58").Value = -9756.000100000001

$ws.Range("H131").Value = 78032.30499999999
$ws.Range("I131").Value = 167118.33
$ws.Range("J131").Value = 51306.5
$ws.Range("K131").Value = 501354.99
$ws.Range("L131").Value = 153919.5
$ws.Range("M131").Value = -496314.99
$ws.Range("N131").Value = -163999.5

$ws.Range("H135").Value = 1751631.4
$ws.Range("I135").Value = 1479710.4
$ws.Range("K135").Value = 13317393.6
$ws.Range("M135").Value = -13314858.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 386.5
$ws.Range("I10").Value = 386.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 386.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -217.5
$ws.Range("N10").ClearContents()

$ws.Range("H12").Value = 251.5
$ws.Range("I12").Value = 251.5
$ws.Range("K12").Value = 251.5
$ws.Range("M12").Value = -111.5

$ws.Range("H132").Value = 22086698
$ws.Range("I132").Value = 23001092
$ws.Range("J132").Value = 21255430
$ws.Range("K132").Value = 69003276
$ws.Range("L132").Value = 63766290
$ws.Range("M132").Value = -69000746
$ws.Range("N132").Value = -63771350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 900
$ws.Range("I17").Value = 900
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -730
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 844427.3
$ws.Range("I132").Value = 318094.03
$ws.Range("J132").Value = 2758366.8
$ws.Range("K132").Value = 954282.0900000001
$ws.Range("L132").Value = 8275100.399999999
$ws.Range("M132").Value = -951752.0900000001
$ws.Range("N132").Value = -8280160.399999999

